$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new D-column value is a plain decimal number that Excel would
# otherwise auto-convert to a numeric type. Force them to remain text (as in
# the original workbook, where every Price cell is stored as a string) by
# setting an explicit text number format before assigning the value.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "29.726.60"
$ws.Range("E2").Value = "  +2.66%  "

# Row 3
$ws.Range("D3").Value = "1.853.15"
$ws.Range("E3").Value = "  +2.00%  "

# Row 4
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "244.84"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").Value = "0.6372"
$ws.Range("E6").Value = "  +4.48%  "

# Row 7
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "0.07530"
$ws.Range("E8").Value = "  +3.26%  "

# Row 9
$ws.Range("D9").Value = "0.2973"
$ws.Range("E9").Value = "  +3.60%  "

# Row 10
$ws.Range("D10").Value = "24.03"
$ws.Range("E10").Value = "  +5.40%  "

# Row 11
$ws.Range("E11").Value = "  +0.50%  "

# Row 12
$ws.Range("D12").Value = "1.839.33"
$ws.Range("E12").Value = "  +1.20%  "

# Row 13
$ws.Range("D13").Value = "5.052"
$ws.Range("E13").Value = "  +2.61%  "

# Row 14
$ws.Range("D14").Value = "0.6896"
$ws.Range("E14").Value = "  +4.82%  "

# Row 15
$ws.Range("D15").Value = "84.73"

# Row 16
$ws.Range("D16").Value = "0.000009596"
$ws.Range("E16").Value = "  +7.91%  "

# Row 17
$ws.Range("D17").Value = "6.059"
$ws.Range("E17").Value = "  +3.72%  "

# Row 18
$ws.Range("D18").Value = "29.706.46"
$ws.Range("E18").Value = "  +2.66%  "

# Row 19
$ws.Range("D19").Value = "2.088.74"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("D20").Value = "239.69"
$ws.Range("E20").Value = "  +1.51%  "

# Row 21
$ws.Range("D21").Value = "12.63"
$ws.Range("E21").Value = "  +1.90%  "

# Row 22
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23
$ws.Range("D23").Value = "7.359"
$ws.Range("E23").Value = "  +3.73%  "

# Row 24
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.12%  "

# Row 25
$ws.Range("D25").Value = "159.58"
$ws.Range("E25").Value = "  +0.76%  "

# Row 26
$ws.Range("D26").Value = "0.1425"
$ws.Range("E26").Value = "  +1.72%  "

# Row 27
$ws.Range("D27").Value = "8.540"
$ws.Range("E27").Value = "  +1.81%  "

# Row 28
$ws.Range("E28").Value = "  +2.19%  "

# Row 29
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  +1.79%  "

# Row 30
$ws.Range("D30").Value = "0.06016"
$ws.Range("E30").Value = "  +7.88%  "

# Row 31
$ws.Range("D31").Value = "1.259"
$ws.Range("E31").Value = "  +4.14%  "

# Row 32
$ws.Range("D32").Value = "4.156"
$ws.Range("E32").Value = "  +2.48%  "

# Row 33
$ws.Range("D33").Value = "4.139"
$ws.Range("E33").Value = "  +1.63%  "

# Row 34
$ws.Range("D34").Value = "1.873"
$ws.Range("E34").Value = "  +3.17%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7333"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.150"
$ws.Range("E36").Value = "  +2.07%  "

# Row 37
$ws.Range("D37").Value = "2.609"
$ws.Range("E37").Value = "  -0.81%  "

# Row 38
$ws.Range("D38").Value = "2.868"
$ws.Range("E38").Value = "  +2.06%  "

# Row 39
$ws.Range("D39").Value = "1.228.53"
$ws.Range("E39").Value = "  +3.09%  "

# Row 40
$ws.Range("E40").Value = "  +1.79%  "

# Row 41
$ws.Range("D41").Value = "6.370"
$ws.Range("E41").Value = "  +0.67%  "

# Row 42
$ws.Range("D42").Value = "0.9165"
$ws.Range("E42").Value = "  +3.16%  "

# Row 43
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("D44").Value = "2.004.11"
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").Value = "102.08"
$ws.Range("E45").Value = "  +1.70%  "

# Row 46
$ws.Range("D46").Value = "66.30"
$ws.Range("E46").Value = "  +3.34%  "

# Row 47
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("E47").Value = "  +2.65%  "

# Row 48
$ws.Range("D48").Value = "0.5079"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49
$ws.Range("D49").Value = "9.305"
$ws.Range("E49").Value = "  +3.04%  "

# Row 50
$ws.Range("D50").Value = "0.4089"
$ws.Range("E50").Value = "  +3.01%  "

# Row 51
$ws.Range("E51").Value = "  +4.00%  "
